$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Narrow columns A:C (previously ~39.29 chars) down to ~34.43 chars wide.
$ws.Columns("A:C").ColumnWidth = 33.71

# Update the forecast figures in column O for a few rows.
$ws.Range("O5").Value = 28.6
$ws.Range("O6").Value = 33.6
$ws.Range("O8").Value = 71.2
$ws.Range("O12").Value = 16.1

# Leave the sheet's selection on M23, matching the saved view state.
$ws.Range("M23").Select() | Out-Null
